{"js": "// \"Lost Legends Errata\" \u2014 the Elementalist class-skill entry that used to\n// read \"Elementalist - Tide\" is renamed to \"Elementalist \u2013 Nature's Embrace\"\n// and the original \"Elementalist - Tide\" entry is kept as its own new\n// sibling bullet directly underneath (same NoSpacing style / ilvl 1, numId 3\n// list formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst OLD_TEXT = \"Elementalist - Tide\";\nconst NEW_TEXT = \"Elementalist \\u2013 Nature\\u2019s Embrace\"; // en dash + curly apostrophe\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === OLD_TEXT) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(`Could not find paragraph with text \"${OLD_TEXT}\"`);\n}\n\n// Insert a sibling paragraph right after the target, inheriting its\n// paragraph formatting (style + numbering), carrying the original text.\ntarget.insertParagraph(OLD_TEXT, \"After\");\n\n// Rename the original paragraph in place to the new errata title.\ntarget.insertText(NEW_TEXT, \"Replace\");\n\nawait context.sync();\n", "ps1": "# \"Lost Legends Errata\" \u2014 the Elementalist class-skill bullet that used to\n# read \"Elementalist - Tide\" is renamed to \"Elementalist \u2013 Nature's Embrace\"\n# and the original \"Elementalist - Tide\" wording is kept as a brand-new\n# sibling bullet directly below it (same NoSpacing style / ilvl 1, numId 3\n# list formatting it always had).\n\n$d = $word.ActiveDocument\n\n$OLD_TEXT = \"Elementalist - Tide\"\n$NEW_TEXT = \"Elementalist \" + [char]0x2013 + \" Nature\" + [char]0x2019 + \"s Embrace\"  # en dash + curly apostrophe\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $OLD_TEXT) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph with text '$OLD_TEXT'\"\n}\n\n# Insert a new sibling paragraph right after the target; Word duplicates the\n# source paragraph's formatting (style + numPr) onto the new paragraph.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = $OLD_TEXT\n\n# Rename the original paragraph in place to the new errata title.\n$target.Range.Text = $NEW_TEXT\n"}
